{"js": "// Office.js (Word JavaScript API) edit script.\n// Applies the Test Log table content changes described by the diff:\n//  - Test 1 row: fill in \"Actual Result\" (replacing the stray _GoBack\n//    bookmark paragraph) and \"Action Required\" cells.\n//  - Test 5 row: clean up the \"wishlist\" spell-check split run in\n//    Description, and fill in \"Actual Result\" / \"Action Required\".\n//  - Test 6 row: clean up the \"wishlist\" spell-check split runs in\n//    Description and Expected Outcome.\n//  - Test 7 row: clean up the \"wishlist\" spell-check split run in\n//    Description.\n//  - Test 8 row: re-insert the _GoBack bookmark in the middle of the\n//    Description text (\"Program shou|ld run without errors\").\n\nconst tables = context.document.body.tables;\ntables.load(\"items\");\nawait context.sync();\n\nconst table = tables.items[0];\n\n// Column indices (0-based): 0=Test #, 1=Description, 2=Expected Outcome,\n// 3=Actual Result, 4=Action Required.\n\n// --- Test 1 (row index 1) ---------------------------------------------\n// \"Actual Result\" cell currently just holds the leftover _GoBack\n// bookmark in an empty paragraph; replacing its body text drops the\n// bookmark and inserts the real text as a single run.\ntable\n  .getCell(1, 3)\n  .body.insertText(\n    \"Text is entered by the by the user into the fields titled ID and Title\",\n    Word.InsertLocation.replace\n  );\n\n// \"Action Required\" cell was empty.\ntable\n  .getCell(1, 4)\n  .body.insertText(\n    \"Ensure when connected to search functions, variables are correct for the required fields\",\n    Word.InsertLocation.replace\n  );\n\n// --- Test 5 (row index 5) ----------------------------------------------\n// Description: merge the spell-check-split \"wishlist\" run back into a\n// single plain run with the same text.\ntable\n  .getCell(5, 1)\n  .body.insertText(\"Save multiple films to a wishlist\", Word.InsertLocation.replace);\n\n// \"Actual Result\" cell was empty.\ntable\n  .getCell(5, 3)\n  .body.insertText(\n    \"Items can be added to wishlist and wishlist is saved, but files are overwritten when a new film is saved.\",\n    Word.InsertLocation.replace\n  );\n\n// \"Action Required\" cell was empty.\ntable\n  .getCell(5, 4)\n  .body.insertText(\n    \"Ensure the file updates to add rather than replace films in the wishlist\",\n    Word.InsertLocation.replace\n  );\n\n// --- Test 6 (row index 6) ----------------------------------------------\ntable\n  .getCell(6, 1)\n  .body.insertText(\n    \"Use the wishlist to quickly access information on saved films\",\n    Word.InsertLocation.replace\n  );\n\ntable\n  .getCell(6, 2)\n  .body.insertText(\n    \"The wishlist should display information about multiple films in full, without requiring another search by the user.\",\n    Word.InsertLocation.replace\n  );\n\n// --- Test 7 (row index 7) ----------------------------------------------\ntable\n  .getCell(7, 1)\n  .body.insertText(\"Delete items from the wishlist\", Word.InsertLocation.replace);\n\n// --- Test 8 (row index 8) ----------------------------------------------\n// Split \"Program should run without errors\" after \"shou\" and drop the\n// _GoBack bookmark back in at that point (Word keeps _GoBack pinned to\n// the most recent edit location).\nconst descCell8 = table.getCell(8, 1);\nconst searchResults = descCell8.body.search(\"shou\", { matchCase: true });\nsearchResults.load(\"text\");\nawait context.sync();\n\nconst insertionPoint = searchResults.items[0].getRange(\"End\");\ninsertionPoint.insertBookmark(\"_GoBack\");\n\nawait context.sync();\n", "ps1": "# Word COM interop edit script.\n# Applies the Test Log table content changes described by the diff:\n#  - Test 1 row: fill in \"Actual Result\" (replacing the stray _GoBack\n#    bookmark paragraph) and \"Action Required\" cells.\n#  - Test 5 row: clean up the \"wishlist\" spell-check split run in\n#    Description, and fill in \"Actual Result\" / \"Action Required\".\n#  - Test 6 row: clean up the \"wishlist\" spell-check split runs in\n#    Description and Expected Outcome.\n#  - Test 7 row: clean up the \"wishlist\" spell-check split run in\n#    Description.\n#  - Test 8 row: re-insert the _GoBack bookmark in the middle of the\n#    Description text (\"Program shou|ld run without errors\").\n\n$doc = $word.ActiveDocument\n$table = $doc.Tables.Item(1)\n\n# Table cells are 1-indexed and include the header row, so data Row N\n# (Test N) lives at Cell(N + 1, col). Columns: 1=Test #, 2=Description,\n# 3=Expected Outcome, 4=Actual Result, 5=Action Required.\n\nfunction Set-CellText($cell, [string]$text) {\n    # Clear the cell's whole range (this also removes any leftover\n    # markup such as proofErr spell-check wrappers or bookmarks) and\n    # then insert the replacement text as a single clean run.\n    $cell.Range.Delete()\n    $cell.Range.InsertBefore($text)\n}\n\n# --- Test 1 (table row 2) ----------------------------------------------\n# \"Actual Result\" cell currently just holds the leftover _GoBack\n# bookmark in an empty paragraph; clearing it drops the bookmark and\n# the new InsertBefore adds the real text as a single run.\nSet-CellText $table.Cell(2, 4) \"Text is entered by the by the user into the fields titled ID and Title\"\n\n# \"Action Required\" cell was empty.\nSet-CellText $table.Cell(2, 5) \"Ensure when connected to search functions, variables are correct for the required fields\"\n\n# --- Test 5 (table row 6) -----------------------------------------------\n# Description: merge the spell-check-split \"wishlist\" run back into a\n# single plain run with the same text.\nSet-CellText $table.Cell(6, 2) \"Save multiple films to a wishlist\"\n\n# \"Actual Result\" cell was empty.\nSet-CellText $table.Cell(6, 4) \"Items can be added to wishlist and wishlist is saved, but files are overwritten when a new film is saved.\"\n\n# \"Action Required\" cell was empty.\nSet-CellText $table.Cell(6, 5) \"Ensure the file updates to add rather than replace films in the wishlist\"\n\n# --- Test 6 (table row 7) -----------------------------------------------\nSet-CellText $table.Cell(7, 2) \"Use the wishlist to quickly access information on saved films\"\nSet-CellText $table.Cell(7, 3) \"The wishlist should display information about multiple films in full, without requiring another search by the user.\"\n\n# --- Test 7 (table row 8) -----------------------------------------------\nSet-CellText $table.Cell(8, 2) \"Delete items from the wishlist\"\n\n# --- Test 8 (table row 9) -----------------------------------------------\n# Split \"Program should run without errors\" after \"shou\" and drop the\n# _GoBack bookmark back in at that point (Word keeps _GoBack pinned to\n# the most recent edit location). Compute the split position from the\n# cell's live Start offset (recalculated after the earlier edits above)\n# plus the length of the \"Program shou\" prefix.\n$descCell8 = $table.Cell(9, 2)\n$prefix = \"Program shou\"\n$splitPos = $descCell8.Range.Start + $prefix.Length\n$bookmarkRange = $doc.Range($splitPos, $splitPos)\n$doc.Bookmarks.Add(\"_GoBack\", $bookmarkRange)\n"}
